# TC04_Canine_Filter_Study-ALL.xlsx - "startup" sheet update
#
# The FilesTab query (row 4, column B "query") is rewritten to drop the
# `File Type` output column and the `Breed` output column (and to stop
# joining on demo.breed), and the active selection on the sheet moves
# from B4 to D4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newFilesQuery = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['COTC007B','NCATS-COP01','GLIOMA01']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS ``File Name``, 
        coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``
"@

$ws.Range("B4").Value = $newFilesQuery

# Selection moves from B4 to D4
$ws.Range("D4").Select()
